$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) from the last existing data row (row 9) down
# to the two new rows (10-11) so the new cells carry the same styles
# (date format on column A, general text style on B:J).
$ws.Range("A9:J9").Copy()
$ws.Range("A10:J11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 10
$ws.Range("A10").Value = 44263.43841596065
$ws.Range("B10").Value = "Session 1 - March 9, Session 2 - March 23, Session 3 - April 6"
$ws.Range("C10").Value = "Novice"
$ws.Range("D10").Value = "Novice"
$ws.Range("E10").Value = "Novice"
$ws.Range("F10").Value = "Novice"
$ws.Range("G10").Value = "Novice"
$ws.Range("H10").Value = "Novice"
$ws.Range("I10").Value = "Novice"
$ws.Range("J10").Value = "Become proficient in using R for small scale projects"

# Row 11
$ws.Range("A11").Value = 44264.416566423606
$ws.Range("B11").Value = "Session 1 - March 9, Session 2 - March 23, Session 3 - April 6"
$ws.Range("C11").Value = "Novice"
$ws.Range("D11").Value = "Basic knowledge or little to none"
$ws.Range("E11").Value = "Novice"
$ws.Range("F11").Value = "Novice"
$ws.Range("G11").Value = "Novice"
$ws.Range("H11").Value = "Basic knowledge or little to none"
$ws.Range("I11").Value = "Novice"
$ws.Range("J11").Value = "To gain more R knowledge so I can better utilize it at work!"
